# Fruta / hortaliza, semanal
# Insert a new weekly record at row 72 (pushing existing rows 72:176 down to 73:177)
# and populate it with the new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(72).Insert()

$ws.Cells.Item(72, 1).Value = 5
$ws.Cells.Item(72, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(72, 3).Value = "Maule"
$ws.Cells.Item(72, 4).Value = 44579
$ws.Cells.Item(72, 5).Value = 7
$ws.Cells.Item(72, 6).Value = "Fruta"
$ws.Cells.Item(72, 7).Value = 100103
$ws.Cells.Item(72, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(72, 9).Value = 100103001
$ws.Cells.Item(72, 10).Value = "Cereza"
$ws.Cells.Item(72, 11).Value = "Lapins"
$ws.Cells.Item(72, 12).Value = "Primera"
$ws.Cells.Item(72, 13).Value = 270
$ws.Cells.Item(72, 14).Value = 4000
$ws.Cells.Item(72, 15).Value = 5000
$ws.Cells.Item(72, 16).Value = 4630
$ws.Cells.Item(72, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(72, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(72, 19).Value = 463
$ws.Cells.Item(72, 20).Value = 10
